$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.041.09'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.646.28'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = "'217.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = "'0.5179"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = "'0.2616"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = "'0.06285"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = "'20.26"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').Value = "'0.07648"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('D12').Value = "'4.561"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('D13').Value = '1.651.30'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '1.873.19'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = "'0.5553"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').Value = '0.0₅8089'
$ws.Range('E16').Value = '  +1.23%  '
$ws.Range('D17').Value = "'64.84"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '26.014.19'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').Value = "'4.582"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').Value = "'10.39"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.08%  '
$ws.Range('D22').Value = "'191.21"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = "'5.886"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.06%  '
$ws.Range('D24').Value = "'1.004"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').Value = "'143.78"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('D26').Value = "'0.1179"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('D27').Value = "'7.156"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = "'15.78"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('D29').Value = "'1.518"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('D30').Value = "'0.05332"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.32%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = "'3.437"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('D33').Value = "'3.326"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').Value = "'1.541"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.46%  '
$ws.Range('D35').Value = "'2.415"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('D36').Value = "'2.777"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('D37').Value = "'0.9377"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').Value = "'0.5566"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').Value = "'0.01570"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').Value = "'5.764"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('D41').Value = "'1.003"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').Value = '1.028.52'
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('D43').Value = "'0.8213"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('D44').Value = "'100.70"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').Value = '1.782.72'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  +6.83%  '
$ws.Range('D47').Value = "'57.03"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = "'1.001"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('D49').Value = "'0.4308"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').Value = "'7.893"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').Value = "'0.05106"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.03%  '
